$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row: "_old" suffix -> "_FV2310", "_new" suffix -> "_FV2404" ---
# Columns A:J carry the "_old" headers, column K is "diff" (unchanged),
# columns L:U carry the "_new" headers.
$oldHeaders = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $colOld = 1 + $i        # A..J
    $colNew = 12 + $i       # L..U
    $ws.Cells.Item(1, $colOld).Value = ($oldHeaders[$i] + "_FV2310")
    $ws.Cells.Item(1, $colNew).Value = ($oldHeaders[$i] + "_FV2404")
}
# Column K ("diff") stays the same.

# --- 2. Turn the used range into a real Excel Table (ListObject) ---
$rng = $ws.Range("A1:U68")
$lo = $ws.ListObjects.Add(1, $rng, [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row (split/freeze pane below row 1) ---
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
